# Operative report template: move the "date" placeholder that used to sit
# alone in F2 into a merged C2:F2 cell (indented with leading spaces so it
# reads like a right-aligned "      {d[i].values[i].date}" label), and make
# the merged cell use a plain Arial 11 font instead of the old bold-ish
# Times New Roman 12 font. Also move the active selection to C2 to match
# where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the date placeholder text currently living in F2 before we touch it.
$dateVal = $ws.Range("F2").Text

# Clear the old standalone F2 value; the other cells in the row (C2:E2)
# were just blank " " filler cells sharing B2's style and get folded into
# the new merged range below.
$ws.Range("F2").Value = ""

# Re-home the date placeholder in C2, offset with a few leading spaces.
$ws.Range("C2").Value = "      " + $dateVal

# Merge C2:F2 into a single cell (this also unifies their formatting).
$ws.Range("C2:F2").Merge()

# Give the merged cell the lighter Arial 11 look instead of the previous
# bold Times New Roman 12 styling (the right-alignment inherited from the
# original F2 style carries through the merge automatically).
$ws.Range("C2:F2").Font.Name = "Arial"
$ws.Range("C2:F2").Font.Size = 11
$ws.Range("C2:F2").Font.Bold = $false

# Move the active selection to the cell that was just edited.
$null = $ws.Range("C2").Select()
